# Fruta / hortaliza, semanal
# Insert this week's new price records (2022-05-06, serial 44687) for
# "Granada" (Wonderfull variety) at the top of the data block, pushing the
# previously-existing rows down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 3 new rows: rows 7-25 shift down to rows 10-28.
$ws.Rows.Item(7).Resize(3).EntireRow.Insert()

$newRows = @(
    @{ K = "Wonderfull"; L = "Especial"; M = 220; N = 21000; O = 21000; P = 21000; Q = "`$/caja 18 kilos granel"; R = "Región de O'Higgins"; S = 1167; T = 18 },
    @{ K = "Wonderfull"; L = "Primera";  M = 250; N = 15000; O = 15000; P = 15000; Q = "`$/caja 18 kilos granel"; R = "Región de O'Higgins"; S = 833;  T = 18 },
    @{ K = "Wonderfull"; L = "Segunda";  M = 280; N = 10000; O = 10000; P = 10000; Q = "`$/caja 18 kilos granel"; R = "Región de O'Higgins"; S = 556;  T = 18 }
)

$r = 7
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = 9
    $ws.Cells.Item($r, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 4).Value = (Get-Date -Year 2022 -Month 5 -Day 6 -Hour 0 -Minute 0 -Second 0)
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100104
    $ws.Cells.Item($r, 8).Value = "Frutos de pepita"
    $ws.Cells.Item($r, 9).Value = 100104001
    $ws.Cells.Item($r, 10).Value = "Granada"
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $r++
}
